$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G4").Value = 0.3552643978832445
$ws.Range("H2:H4").Value = 0.9990000000000001
